$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D:D").Insert()

$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D7").Value = 43465
Write-Output $ws.Range("D7").Text
Write-Output $ws.Range("D11").Text
$ws.Range("D16").Value = $null
Write-Output "D16 ok"
